$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing last column (K) into the new column L for each row
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)

$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("L3").Value = 2022

$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 370

$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 137

$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value = 314

$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value = 121

$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L8").NumberFormat = "#,##0"
$ws.Range("L8").Value = 50

$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L9").Value = 16

$ws.Range("L2").Select() | Out-Null
